$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Khloe Bright"
$ws.Range("B4").Value = 18
$ws.Range("C4").Value = "Female"
$ws.Range("D4").Value = "Insects"
$ws.Range("E4").Value = "Female Bestfriend"

$ws.Range("A5").Value = "Percy Thomson"
$ws.Range("B5").Value = 17
$ws.Range("C5").Value = "Male"
$ws.Range("D5").Value = "Mythicals(Elves)"
$ws.Range("E5").Value = "Comical Partner"
